# The commit swaps the presentation's two embedded themes: the theme that is
# actually wired up to the slide master / presentation (holding the
# "Integral" green colour scheme) is replaced with the colours of the
# "Office Theme" default scheme (which, before the edit, sat unused in the
# companion theme part).
#
# PowerPoint exposes the live theme of the active design through
# SlideMaster.Theme.ThemeColorScheme - each of the 12 items corresponds,
# in order, to dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink of the
# <a:clrScheme> element. Setting .RGB on each item rewrites those colours.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target values: the stock "Office Theme" colour scheme.
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i]
    # PowerPoint's RGB long uses 0x00BBGGRR byte order.
    $r = ($hex -band 0xFF0000) -shr 16
    $g = ($hex -band 0x00FF00) -shr 8
    $b = ($hex -band 0x0000FF)
    $bgr = $r + ($g * 256) + ($b * 65536)

    $colorScheme.Item($i + 1).RGB = $bgr
}
